$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, centered, bordered) from H1 onto the two new
# header cells before writing their text, so I1/J1 pick up style index 1
# just like the other header cells (B1:H1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows for the two new columns.
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 6

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 4

$ws.Range("I4").Value = 4
$ws.Range("J4").Value = 6

$ws.Range("I5").Value = 4
$ws.Range("J5").Value = 6

$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 2
